$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.016.96"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "3.926.74"
$ws.Range("E3").Value = "  +4.85%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'604.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'165.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").Value = "3.924.72"
$ws.Range("E7").Value = "  +4.87%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("D10").Value = "'0.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.12%  "
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "'37.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "4.588.44"
$ws.Range("E15").Value = "  +5.00%  "
$ws.Range("D16").Value = "3.948.49"
$ws.Range("E16").Value = "  +5.42%  "
$ws.Range("D17").Value = "69.108.79"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "'7.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").Value = "'17.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.83%  "
$ws.Range("D21").Value = "'11.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("D22").Value = "'489.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  +12.30%  "
$ws.Range("D25").Value = "'84.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "'2.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("D28").Value = "'10.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").Value = "4.079.02"
$ws.Range("E31").Value = "  +4.93%  "
$ws.Range("D32").Value = "'7.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.60%  "
$ws.Range("D33").Value = "'2.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.74%  "
$ws.Range("D34").Value = "'32.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("D35").Value = "3.875.23"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("E37").Value = "  +2.87%  "
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("D39").Value = "'5.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("D42").Value = "'3.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.47%  "
$ws.Range("D43").Value = "'441.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.11%  "
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D47").Value = "'8.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("D48").Value = "'27.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +15.42%  "
$ws.Range("D49").Value = "2.853.69"
$ws.Range("E49").Value = "  +2.37%  "
$ws.Range("D50").Value = "'142.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("E51").Value = "  +2.12%  "
